$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("0ce5dd49", "During work - study hours"),
    @("2c1001cb", "During work - study hours"),
    @("37cc37bf", "During work - study hours"),
    @("43faa0b9", "During work - study hours"),
    @("4abe3e88", "During work - study hours"),
    @("50164f59", "During work - study hours"),
    @("5cf70f79", "During work - study hours"),
    @("5da96769", "During work - study hours"),
    @("6ca3e2f6", "During work - study hours"),
    @("790a4fcb", "During work - study hours"),
    @("802cc63a", "During work - study hours"),
    @("85c3ea4d", "Anytime throughout the day"),
    @("942dfafb", "During work - study hours"),
    @("9bc6ba8c", "During work - study hours"),
    @("a2d65af2", "During work - study hours"),
    @("a46f1771", "During work - study hours"),
    @("ad58f9da", "During work - study hours"),
    @("c7d9a301", "Anytime throughout the day"),
    @("ce8732ff", "During work - study hours"),
    @("d6f1d567", "During work - study hours"),
    @("da9326c9", "During work - study hours"),
    @("e09ca7bf", "During work - study hours"),
    @("ef53a641", "During work - study hours")
)

$timestamp = 45854.65031743384
$question = "q04_most_common_time"

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $question
    $ws.Cells.Item($row, 4).Value = $timestamp
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $row++
}
